$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Add 4 new rows to the "Tableau4" table (rows 78-81) ---
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Copy cell formatting (styles/borders) from existing template rows ---
# Row 78 reuses the "banded" formatting seen on row 21
$ws.Range("A21:H21").Copy() | Out-Null
$ws.Range("A78:H78").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Rows 79-81 reuse the formatting of row 75 (the usual last-row style)
$ws.Range("A75:H75").Copy() | Out-Null
$ws.Range("A79:H81").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 78 values ---
$ws.Cells.Item(78, 1).Value = 44643
$ws.Cells.Item(78, 2).Value = 0.4152777777777778
$ws.Cells.Item(78, 3).Value = 0.42152777777777778
$ws.Cells.Item(78, 4).Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Cells.Item(78, 5).Value = "CPNV"
$ws.Cells.Item(78, 6).Value = "Commenter"
$ws.Cells.Item(78, 7).Value = "J'ai organizer un peu mon code et je l'ai comenté"

# --- Row 80 Description is typed before row 79's, matching the source edit history ---
$ws.Cells.Item(80, 7).Value = "J'ai crée une fonction qui choisi une carte aléatoire (entre 1 et 5)"

# --- Row 79 values ---
$ws.Cells.Item(79, 1).Value = 44643
$ws.Cells.Item(79, 2).Value = 0.42222222222222222
$ws.Cells.Item(79, 3).Value = 0.42708333333333331
$ws.Cells.Item(79, 4).Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Cells.Item(79, 5).Value = "CPNV"
$ws.Cells.Item(79, 6).Value = "Créer les 5 grilles"
$ws.Cells.Item(79, 7).Value = "J'ai crée les 5 grilles"

# --- Row 80 remaining values ---
$ws.Cells.Item(80, 1).Value = 44643
$ws.Cells.Item(80, 2).Value = 0.42708333333333331
$ws.Cells.Item(80, 3).Value = 0.46388888888888885
$ws.Cells.Item(80, 4).Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Cells.Item(80, 5).Value = "CPNV"
$ws.Cells.Item(80, 6).Value = "Fonction randomGrid"

# --- Row 81 values ---
$ws.Cells.Item(81, 1).Value = 44643
$ws.Cells.Item(81, 2).Value = 0.46597222222222223
$ws.Cells.Item(81, 3).Value = 0.4694444444444445
$ws.Cells.Item(81, 4).Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Cells.Item(81, 5).Value = "CPNV"
$ws.Cells.Item(81, 6).Value = "Commenter"
$ws.Cells.Item(81, 7).Value = "J'ai commenté et testé une dernière fois ma fonction randomGrid"

# --- H column stays blank on all 4 new rows (already blank after paste) ---

# --- Update the conditional-formatting range that excludes a stale block of rows ---
$cfRng = $ws.Range("D233:D1048576,D1:D206")
$cfs = $cfRng.FormatConditions
for ($i = 1; $i -le $cfs.Count; $i++) {
    $cf = $cfs.Item($i)
    if ($cf.Operator -eq 7) {
        $cf.ModifyAppliesToRange($ws.Range("D236:D1048576,D1:D209"))
    }
}

# --- Update the view: scroll position and active selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
$ws.Range("G86").Select()
